$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the helper table that had been added in columns B:O (rows 3-5)
# and the raw per-observation values spilled across columns D:O (rows 4-5).
# Columns D4:F4 and C5 keep their existing "0.0" number format but become
# blank; everything else in B3:O5 is cleared entirely.
$ws.Range("B3:O5").ClearContents()

# The extra font style (white, Arial Unicode MS, vertical-center) that had
# been applied to column A is no longer used - restore the default style.
$ws.Range("A1:A27").ClearFormats()

# Update the active selection to match the reduced used range.
$ws.Range("F10").Select()
